$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-26 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-27 Saturday", 2)

$d.Content.Find.Execute("25×11=275", $true, $false, $false, $false, $false, $true, 1, $false, "40×37=1480", 2)
$d.Content.Find.Execute("55×24=1320", $true, $false, $false, $false, $false, $true, 1, $false, "37×24=888", 2)
$d.Content.Find.Execute("99×11=1089", $true, $false, $false, $false, $false, $true, 1, $false, "61×27=1647", 2)
$d.Content.Find.Execute("59×87=5133", $true, $false, $false, $false, $false, $true, 1, $false, "92×31=2852", 2)
$d.Content.Find.Execute("78×73=5694", $true, $false, $false, $false, $false, $true, 1, $false, "56×84=4704", 2)

$d.Content.Find.Execute("39×78=3042", $true, $false, $false, $false, $false, $true, 1, $false, "35×28=980", 2)
$d.Content.Find.Execute("17×49=833", $true, $false, $false, $false, $false, $true, 1, $false, "49×82=4018", 2)
$d.Content.Find.Execute("15×58=870", $true, $false, $false, $false, $false, $true, 1, $false, "24×13=312", 2)
$d.Content.Find.Execute("37×74=2738", $true, $false, $false, $false, $false, $true, 1, $false, "79×14=1106", 2)
$d.Content.Find.Execute("27×43=1161", $true, $false, $false, $false, $false, $true, 1, $false, "66×58=3828", 2)

$d.Content.Find.Execute("35×46=1610", $true, $false, $false, $false, $false, $true, 1, $false, "24×89=2136", 2)
$d.Content.Find.Execute("50×41=2050", $true, $false, $false, $false, $false, $true, 1, $false, "17×51=867", 2)
$d.Content.Find.Execute("22×55=1210", $true, $false, $false, $false, $false, $true, 1, $false, "67×59=3953", 2)
$d.Content.Find.Execute("66×66=4356", $true, $false, $false, $false, $false, $true, 1, $false, "90×23=2070", 2)
$d.Content.Find.Execute("59×34=2006", $true, $false, $false, $false, $false, $true, 1, $false, "12×15=180", 2)

$d.Content.Find.Execute("44×13=572", $true, $false, $false, $false, $false, $true, 1, $false, "16×91=1456", 2)
$d.Content.Find.Execute("46×81=3726", $true, $false, $false, $false, $false, $true, 1, $false, "54×96=5184", 2)
$d.Content.Find.Execute("43×45=1935", $true, $false, $false, $false, $false, $true, 1, $false, "43×85=3655", 2)
$d.Content.Find.Execute("24×63=1512", $true, $false, $false, $false, $false, $true, 1, $false, "43×54=2322", 2)
$d.Content.Find.Execute("68×54=3672", $true, $false, $false, $false, $false, $true, 1, $false, "30×44=1320", 2)

$d.Content.Find.Execute("66×45=2970", $true, $false, $false, $false, $false, $true, 1, $false, "88×75=6600", 2)
$d.Content.Find.Execute("54×53=2862", $true, $false, $false, $false, $false, $true, 1, $false, "95×44=4180", 2)
$d.Content.Find.Execute("39×41=1599", $true, $false, $false, $false, $false, $true, 1, $false, "77×94=7238", 2)
$d.Content.Find.Execute("46×73=3358", $true, $false, $false, $false, $false, $true, 1, $false, "18×11=198", 2)
$d.Content.Find.Execute("62×15=930", $true, $false, $false, $false, $false, $true, 1, $false, "76×37=2812", 2)
